$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column CM (index 91): give it the same width as the other data columns ---
$ws.Columns(91).ColumnWidth = $ws.Columns(90).ColumnWidth

# --- Template cells already carrying the three cell styles used throughout the sheet ---
# (s="1" plain / s="2" yellow "low value" fill / s="3" blue "mid value" fill)
$styleN = $ws.Range("A2")
$styleL = $ws.Range("D2")
$styleM = $ws.Range("N2")

# --- CM1: new header date, entered as literal text like the rest of row 1 ---
$cm1 = $ws.Range("CM1")
$cm1.NumberFormat = "@"
$cm1.Value = "2024/12/08"
$ws.Range("CL1").Copy()
$cm1.PasteSpecial(-4122)

# --- CM2:CM53: the 2024/12/08 data column, one row per existing data row ---
$data = @(
    @{ Row = 2; Value = 142.7; Style = "N" },
    @{ Row = 3; Value = 170.7; Style = "N" },
    @{ Row = 4; Value = 164.3; Style = "N" },
    @{ Row = 5; Value = 145.1; Style = "N" },
    @{ Row = 6; Value = 165.1; Style = "N" },
    @{ Row = 7; Value = 95.90000000000001; Style = "L" },
    @{ Row = 8; Value = 177.2; Style = "N" },
    @{ Row = 9; Value = 146.8; Style = "N" },
    @{ Row = 10; Value = 151.5; Style = "N" },
    @{ Row = 11; Value = 141.2; Style = "N" },
    @{ Row = 12; Value = 147.6; Style = "N" },
    @{ Row = 13; Value = 176.8; Style = "N" },
    @{ Row = 14; Value = 143.6; Style = "N" },
    @{ Row = 15; Value = 166.2; Style = "N" },
    @{ Row = 16; Value = 132; Style = "M" },
    @{ Row = 17; Value = 225.8; Style = "N" },
    @{ Row = 18; Value = 146.5; Style = "N" },
    @{ Row = 19; Value = 121.7; Style = "L" },
    @{ Row = 20; Value = 146.9; Style = "N" },
    @{ Row = 21; Value = 164; Style = "N" },
    @{ Row = 22; Value = 180.4; Style = "N" },
    @{ Row = 23; Value = 140.2; Style = "N" },
    @{ Row = 24; Value = 195.5; Style = "N" },
    @{ Row = 25; Value = 123.4; Style = "L" },
    @{ Row = 26; Value = 126.1; Style = "M" },
    @{ Row = 27; Value = 300.3; Style = "N" },
    @{ Row = 28; Value = 168.8; Style = "N" },
    @{ Row = 29; Value = 173.3; Style = "N" },
    @{ Row = 30; Value = 144.7; Style = "N" },
    @{ Row = 31; Value = 186.9; Style = "N" },
    @{ Row = 32; Value = 136.8; Style = "M" },
    @{ Row = 33; Value = 194.4; Style = "N" },
    @{ Row = 34; Value = 175.2; Style = "N" },
    @{ Row = 35; Value = 140.3; Style = "N" },
    @{ Row = 36; Value = 174.2; Style = "N" },
    @{ Row = 37; Value = 174.5; Style = "N" },
    @{ Row = 38; Value = 145.7; Style = "N" },
    @{ Row = 39; Value = 198.3; Style = "N" },
    @{ Row = 40; Value = 187.7; Style = "N" },
    @{ Row = 41; Value = 116.5; Style = "L" },
    @{ Row = 42; Value = 179; Style = "N" },
    @{ Row = 43; Value = 158.5; Style = "N" },
    @{ Row = 44; Value = 132.3; Style = "M" },
    @{ Row = 45; Value = 409.8; Style = "N" },
    @{ Row = 46; Value = 152.5; Style = "N" },
    @{ Row = 47; Value = 152.3; Style = "N" },
    @{ Row = 48; Value = 197.4; Style = "N" },
    @{ Row = 49; Value = 158.4; Style = "N" },
    @{ Row = 50; Value = 209.8; Style = "N" },
    @{ Row = 51; Value = 127.3; Style = "M" },
    @{ Row = 52; Value = 160.5; Style = "N" },
    @{ Row = 53; Value = 154.1; Style = "N" }
)

foreach ($row in $data) {
    $dst = $ws.Cells.Item($row.Row, 91)
    $dst.Value = $row.Value
    switch ($row.Style) {
        "N" { $styleN.Copy() }
        "L" { $styleL.Copy() }
        "M" { $styleM.Copy() }
    }
    $dst.PasteSpecial(-4122)
}
